$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 348
$ws.Range("I53").Value = 410.44446
$ws.Range("J53").Value = 291.8
$ws.Range("K53").Value = 410.44446
$ws.Range("L53").Value = 291.8
$ws.Range("M53").Value = 226.55554
$ws.Range("N53").Value = -1565.8
$ws.Range("H118").Value = 999.6667
$ws.Range("J118").Value = 999
$ws.Range("L118").Value = 2997
$ws.Range("N118").Value = -6311
$ws.Range("H133").Value = 125000
$ws.Range("J133").Value = 125000
$ws.Range("L133").Value = 125000
$ws.Range("N133").Value = -135120
$ws.Range("H136").Value = 140220.33
$ws.Range("J136").Value = 140220.33
$ws.Range("L136").Value = 140220.33
$ws.Range("N136").Value = -150420.33
$ws.Range("H137").Value = 1416.8948
$ws.Range("I137").Value = 1171.8823
$ws.Range("K137").Value = 3515.6469
$ws.Range("M137").Value = -965.6468999999997
$ws.Range("H138").Value = 2920.25
$ws.Range("I138").Value = 1931.3
$ws.Range("J138").Value = 4568.5
$ws.Range("K138").Value = 5793.9
$ws.Range("L138").Value = 13705.5
$ws.Range("M138").Value = -653.8999999999996
$ws.Range("N138").Value = -23985.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4997.5054
$ws.Range("I32").Value = 4723.5933
$ws.Range("J32").Value = 9151.833000000001
$ws.Range("K32").Value = 4723.5933
$ws.Range("L32").Value = 9151.833000000001
$ws.Range("M32").Value = -4436.5933
$ws.Range("N32").Value = -9725.833000000001
$ws.Range("H61").Value = 5042.5557
$ws.Range("I61").Value = 10791.75
$ws.Range("J61").Value = 3399.9285
$ws.Range("K61").Value = 10791.75
$ws.Range("L61").Value = 3399.9285
$ws.Range("M61").Value = -10579.75
$ws.Range("N61").Value = -3823.9285
$ws.Range("H74").Value = 2804
$ws.Range("I74").Value = 2765.2
$ws.Range("K74").Value = 2765.2
$ws.Range("M74").Value = -1891.2
$ws.Range("H77").Value = 2804
$ws.Range("I77").Value = 2765.2
$ws.Range("K77").Value = 13826
$ws.Range("M77").Value = -9458
$ws.Range("H132").Value = 2647.4783
$ws.Range("I132").Value = 2519.6
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 7558.799999999999
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -5028.799999999999
$ws.Range("N132").Value = -15560
$ws.Range("H133").Value = 106753.336
$ws.Range("J133").Value = 106753.336
$ws.Range("L133").Value = 106753.336
$ws.Range("N133").Value = -111813.336
$ws.Range("H136").Value = 5042.5557
$ws.Range("I136").Value = 10791.75
$ws.Range("J136").Value = 3399.9285
$ws.Range("K136").Value = 32375.25
$ws.Range("L136").Value = 10199.7855
$ws.Range("M136").Value = -29825.25
$ws.Range("N136").Value = -15299.7855
$ws.Range("H139").Value = 137399.2
$ws.Range("J139").Value = 137399.2
$ws.Range("L139").Value = 137399.2
$ws.Range("N139").Value = -147679.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 13483
$ws.Range("J81").Value = 13483
$ws.Range("L81").Value = 13483
$ws.Range("N81").Value = -15605
$ws.Range("H84").Value = 13483
$ws.Range("J84").Value = 13483
$ws.Range("L84").Value = 40449
$ws.Range("N84").Value = -51057
$ws.Range("H133").Value = 329950
$ws.Range("J133").Value = 329950
$ws.Range("L133").Value = 329950
$ws.Range("N133").Value = -340070
$ws.Range("H134").Value = 3987.6667
$ws.Range("I134").Value = 3987.6667
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11963.0001
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -9428.000100000001
$ws.Range("H135").Value = 210332.83
$ws.Range("J135").Value = 210332.83
$ws.Range("L135").Value = 210332.83
$ws.Range("N135").Value = -220472.83
$ws.Range("H137").Value = 98832.664
$ws.Range("J137").Value = 98832.664
$ws.Range("L137").Value = 98832.664
$ws.Range("N137").Value = -109032.664
$ws.Range("H138").Value = 68650.336
$ws.Range("J138").Value = 82999
$ws.Range("L138").Value = 82999
$ws.Range("N138").Value = -93279
$ws.Range("H140").Value = 38239.832
$ws.Range("J140").Value = 38239.832
$ws.Range("L140").Value = 38239.832
$ws.Range("N140").Value = -48599.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1980.6923
$ws.Range("I31").Value = 1647.3684
$ws.Range("J31").Value = 2885.4285
$ws.Range("K31").Value = 1647.3684
$ws.Range("L31").Value = 2885.4285
$ws.Range("M31").Value = -1352.3684
$ws.Range("N31").Value = -3475.4285
$ws.Range("H34").Value = 1980.6923
$ws.Range("I34").Value = 1647.3684
$ws.Range("J34").Value = 2885.4285
$ws.Range("K34").Value = 1647.3684
$ws.Range("L34").Value = 2885.4285
$ws.Range("M34").Value = -1445.3684
$ws.Range("N34").Value = -3289.4285
$ws.Range("H58").Value = 2376.2144
$ws.Range("I58").Value = 2328.2307
$ws.Range("K58").Value = 2328.2307
$ws.Range("M58").Value = -2125.2307
$ws.Range("H132").Value = 6551.972
$ws.Range("I132").Value = 2135.963
$ws.Range("J132").Value = 19800
$ws.Range("K132").Value = 6407.889000000001
$ws.Range("L132").Value = 59400
$ws.Range("M132").Value = -3877.889000000001
$ws.Range("N132").Value = -64460
$ws.Range("H133").Value = 57519.137
$ws.Range("J133").Value = 59071.05
$ws.Range("L133").Value = 59071.05
$ws.Range("N133").Value = -64131.05
$ws.Range("H134").Value = 3292.2258
$ws.Range("I134").Value = 3112.7637
$ws.Range("J134").Value = 4702.2856
$ws.Range("K134").Value = 9338.2911
$ws.Range("L134").Value = 14106.8568
$ws.Range("M134").Value = -6803.2911
$ws.Range("N134").Value = -19176.8568
$ws.Range("H135").Value = 86682.836
$ws.Range("J135").Value = 86682.836
$ws.Range("L135").Value = 86682.836
$ws.Range("N135").Value = -96822.836
$ws.Range("H136").Value = 2376.2144
$ws.Range("I136").Value = 2328.2307
$ws.Range("K136").Value = 6984.6921
$ws.Range("M136").Value = -4434.6921
$ws.Range("H138").Value = 78505.82000000001
$ws.Range("J138").Value = 78505.82000000001
$ws.Range("L138").Value = 78505.82000000001
$ws.Range("N138").Value = -88785.82000000001
$ws.Range("H140").Value = 105354.25
$ws.Range("J140").Value = 99999.5
$ws.Range("L140").Value = 99999.5
$ws.Range("N140").Value = -110359.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 234290.75
$ws.Range("J134").Value = 234290.75
$ws.Range("L134").Value = 702872.25
$ws.Range("N134").Value = -707942.25
$ws.Range("H135").Value = 145177
$ws.Range("J135").Value = 145177
$ws.Range("L135").Value = 145177
$ws.Range("N135").Value = -155317
$ws.Range("H139").Value = 124990
$ws.Range("J139").Value = 124990
$ws.Range("L139").Value = 124990
$ws.Range("N139").Value = -135270
$ws.Range("H140").Value = 67374.414
$ws.Range("J140").Value = 67374.414
$ws.Range("L140").Value = 67374.414
$ws.Range("N140").Value = -77734.414

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 114872.1
$ws.Range("J133").Value = 114872.1
$ws.Range("L133").Value = 114872.1
$ws.Range("N133").Value = -119932.1
$ws.Range("H135").Value = 87923.25
$ws.Range("J135").Value = 87923.25
$ws.Range("L135").Value = 87923.25
$ws.Range("N135").Value = -98063.25
$ws.Range("H136").Value = 7589.968
$ws.Range("I136").Value = 10099.294
$ws.Range("J136").Value = 4542.9287
$ws.Range("K136").Value = 30297.882
$ws.Range("L136").Value = 13628.7861
$ws.Range("M136").Value = -27747.882
$ws.Range("N136").Value = -18728.7861
$ws.Range("H140").Value = 108991
$ws.Range("J140").Value = 108991
$ws.Range("L140").Value = 108991
$ws.Range("N140").Value = -119351

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 80624.5
$ws.Range("J46").Value = 80624.5
$ws.Range("L46").Value = 80624.5
$ws.Range("N46").Value = -81086.5
$ws.Range("H132").Value = 3001.087
$ws.Range("I132").Value = 2810.6428
$ws.Range("J132").Value = 5000.75
$ws.Range("K132").Value = 8431.928400000001
$ws.Range("L132").Value = 15002.25
$ws.Range("M132").Value = -5901.928400000001
$ws.Range("N132").Value = -20062.25
$ws.Range("H134").Value = 80624.5
$ws.Range("J134").Value = 80624.5
$ws.Range("L134").Value = 241873.5
$ws.Range("N134").Value = -246943.5
$ws.Range("H135").Value = 50534.668
$ws.Range("J135").Value = 50534.668
$ws.Range("L135").Value = 50534.668
$ws.Range("N135").Value = -60674.668
$ws.Range("H137").Value = 50950.375
$ws.Range("J137").Value = 50950.375
$ws.Range("L137").Value = 50950.375
$ws.Range("N137").Value = -61150.375
$ws.Range("H138").Value = 114999.5
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("H140").Value = 60965.168
$ws.Range("J140").Value = 60965.168
$ws.Range("L140").Value = 60965.168
$ws.Range("N140").Value = -71325.16800000001
